$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new value map, derived from the canonical OOXML diff of the workbook.
# Every cell in the source file is a plain text/inline string (t="inlineStr"),
# including price cells that merely look numeric (e.g. "0.999", "7.61"). Plain
# ".Value = ..." assignment is fine for text that Excel cannot parse as a number
# (URLs, names, "%"-suffixed deltas, multi-dot prices like "26.415.84"), but for
# the handful of purely-numeric-looking price strings we must force the cell to
# "@" (text) format first so Excel keeps storing/display them as text instead of
# silently converting them into real numbers.

$plainUpdates = [ordered]@{
    "D2" = "26.415.84"
    "E2" = "  +0.59%  "
    "D3" = "1.608.26"
    "E3" = "  +0.93%  "
    "E4" = "  -0.10%  "
    "E5" = "  -0.75%  "
    "E6" = "  -0.90%  "
    "E7" = "  -0.09%  "
    "B8" = "Cardano"
    "C8" = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
    "E8" = "  -0.66%  "
    "B9" = "Dogecoin"
    "C9" = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
    "E9" = "  -0.31%  "
    "E10" = "  +1.36%  "
    "E11" = "  -0.51%  "
    "D12" = "1.833.72"
    "E12" = "  +0.86%  "
    "D13" = "1.608.81"
    "E13" = "  +0.79%  "
    "E14" = "  -0.34%  "
    "E15" = "  -0.45%  "
    "E16" = "  -0.92%  "
    "D17" = "26.405.90"
    "E17" = "  +0.55%  "
    "E18" = "  +8.11%  "
    "E19" = "  -0.35%  "
    "E20" = "  +3.03%  "
    "E21" = "  -0.14%  "
    "E22" = "  -0.57%  "
    "E23" = "  +4.84%  "
    "E24" = "  -0.83%  "
    "E25" = "  +1.17%  "
    "E26" = "  -0.07%  "
    "E27" = "  -0.12%  "
    "E28" = "  +0.58%  "
    "E29" = "  +1.85%  "
    "E30" = "  +0.73%  "
    "E31" = "  -0.69%  "
    "D32" = "1.492.75"
    "E32" = "  +5.19%  "
    "E33" = "  +0.80%  "
    "E34" = "  -1.35%  "
    "E35" = "  -0.49%  "
    "E36" = "  +0.87%  "
    "E37" = "  -3.01%  "
    "E39" = "  -0.18%  "
    "E40" = "  +0.01%  "
    "E41" = "  -0.04%  "
    "E43" = "  -5.47%  "
    "D44" = "1.746.76"
    "E44" = "  +0.93%  "
    "E45" = "  -0.02%  "
    "E46" = "  -0.19%  "
    "E47" = "  +2.70%  "
    "B48" = "BabyDogeCoin"
    "C48" = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
    "D48" = "0.0₆0104"
    "E48" = "  -1.61%  "
    "B49" = "RenderToken"
    "C49" = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
    "E49" = "  -0.99%  "
    "B50" = "Cronos"
    "C50" = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
    "E50" = "  -1.41%  "
    "B51" = "Algorand"
    "C51" = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
    "E51" = "  +0.39%  "
}

$textUpdates = [ordered]@{
    "D5" = "211.71"
    "D8" = "0.244"
    "D9" = "0.0605"
    "D10" = "19.23"
    "D11" = "0.0849"
    "D16" = "63.39"
    "D18" = "232.07"
    "D20" = "7.61"
    "D21" = "0.999"
    "D23" = "2.20"
    "D24" = "8.98"
    "D25" = "146.66"
    "D29" = "15.42"
    "D30" = "0.0494"
    "D38" = "0.0165"
    "D39" = "0.822"
    "D40" = "5.80"
    "D43" = "0.928"
    "D45" = "0.762"
    "D46" = "60.87"
    "D47" = "89.50"
    "D49" = "1.48"
    "D50" = "0.0501"
    "D51" = "0.0959"
}

foreach ($ref in $plainUpdates.Keys) {
    $ws.Range($ref).Value = $plainUpdates[$ref]
}

foreach ($ref in $textUpdates.Keys) {
    $cell = $ws.Range($ref)
    $cell.Style = "Normal"
    $cell.NumberFormat = "@"
    $cell.Value = $textUpdates[$ref]
}
